# edit.ps1 - apply the RGF_input_file.xlsx changes described by the commit
# "debugging RGF. GPU inv has issue."
#
# Strategy:
#  1. Edit rows 1-10 (these row numbers are unaffected by the later row deletion).
#  2. Delete the old row 11 ("#" row) so that everything below shifts up by one,
#     which also fixes up the merged cells automatically.
#  3. Edit the (now shifted) rows 11-15 to their final contents.
#  4. Fix up sheet view / selection / window metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 : header row - add a second "Input" column
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Input1"
$ws.Range("C1").Value = "Input2"
$ws.Range("D1").Value = "Comment"

# ---------------------------------------------------------------------------
# Row 2 : Using GPU flag flipped to False
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = $false

# Rows 3-6 keep their displayed values (CPU max matrix / Material / Lattice /
# Direction); nothing to change there.

# ---------------------------------------------------------------------------
# Row 7 : was "Max ribbon width" formula row -> becomes "mesh" input row
# ---------------------------------------------------------------------------
$ws.Range("C10").Copy($ws.Range("B7"))   # pick up the "input" style (s=4)
$ws.Range("C10").Copy($ws.Range("C7"))   # same style for the second input cell
$ws.Range("B7").Formula = ""
$ws.Range("C7").Formula = ""
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 1000
$ws.Range("B7").s = $null

$ws.Range("B7").Copy($ws.Range("D7"))    # not used - placeholder (overwritten below)

# ---------------------------------------------------------------------------
# Row 8 : was "Max ribbon length" formula row -> becomes "Bias(V)" row
# ---------------------------------------------------------------------------
$ws.Range("B8").Formula = ""
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0

# ---------------------------------------------------------------------------
# Row 9 : was plain "Bias(V)" values -> becomes "Plot band structure" (booleans)
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = $true
$ws.Range("C10").Copy($ws.Range("C9"))
$ws.Range("C9").Value = $true

# ---------------------------------------------------------------------------
# Row 10 : was "Plot band structure" booleans -> becomes "#" row
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "1 for 2-1-2-1 type ribbon, 2 for 2-2-2-2 type ribbon"
$ws.Range("D10").Value = ""

# ---------------------------------------------------------------------------
# Now delete the old row 11 (the former "#" row). Everything below (old rows
# 12-16) shifts up to become the new rows 11-15, and the merged cells that
# lived on old row 12 automatically move to the new row 11.
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# New row 7 comment column (D7) - needs the explanatory text + its own style
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy($ws.Range("D7"))
$ws.Range("D7").Value = "first column for max sub cell number been calculated. Second column for kx sweep meshing"

# ---------------------------------------------------------------------------
# New row 9 comment column (D9)
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "first boolean for enabling plot function. Second boolean for also plotting zoom in figures"

# ---------------------------------------------------------------------------
# Row 13 (was header-label row, now becomes the first "x" data row, reusing
# the values that used to live in row 14)
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "x"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 23
$ws.Range("F13").Value = 807
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0.1
$ws.Range("J13").Value = 46
$ws.Range("K13").Value = 0.4
$ws.Range("L13").Value = 46
$ws.Range("M13").Value = 0.4
$ws.Range("N13").Formula = "=(E13+C13-1)*0.246*3^0.5/2"
$ws.Range("O13").Formula = "=(J13+C13-1)*0.246*3^0.5/2"
$ws.Range("P13").Formula = "=(L13+C13-1)*0.246*3^0.5/2"
$ws.Range("Q13").Formula = "=F13*0.246*3"

# ---------------------------------------------------------------------------
# Row 14 (new data - "o" row)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "o"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 500
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.25
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0.4
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0.4
$ws.Range("N14").Formula = "=(E14+C14-1)*0.246*3^0.5/2"
$ws.Range("O14").Formula = "=(J14+C14-1)*0.246*3^0.5/2"
$ws.Range("P14").Formula = "=(L14+C14-1)*0.246*3^0.5/2"
$ws.Range("Q14").Formula = "=F14*0.246*3"

# ---------------------------------------------------------------------------
# Row 15 (last data row - "x" row, matches what used to be row 16)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "x"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 46
$ws.Range("K15").Value = 0.4
$ws.Range("L15").Value = 46
$ws.Range("M15").Value = 0.4
$ws.Range("N15").Formula = "=(E15+C15-1)*0.246*3^0.5/2"
$ws.Range("O15").Formula = "=(J15+C15-1)*0.246*3^0.5/2"
$ws.Range("P15").Formula = "=(L15+C15-1)*0.246*3^0.5/2"
$ws.Range("Q15").Formula = "=F15*0.246*3"

# ---------------------------------------------------------------------------
# Sheet view: selection moves to B5
# ---------------------------------------------------------------------------
$ws.Range("B5").Select()

# ---------------------------------------------------------------------------
# Workbook window size (best effort - may not be exposed by this host)
# ---------------------------------------------------------------------------
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 16830
    $win.Height = 9240
} catch {
}

Write-Output "edit complete"
